$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New test case row (sc13) added after row 17
$ws.Range("A18").Value = "sc13"
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 43
$ws.Range("E18").Value = "A(5803) Single hierarchy conditions, two policies. one extra condition on one policy"
$ws.Range("F18").Value = "input files"
$ws.Range("G18").Value = "no"
$ws.Range("H18").Value = "to do"

# Column E (Description) widened to fit the new, longer text
$ws.Columns.Item(5).ColumnWidth = 68.8

# Selection left on G10 after the edit
$ws.Range("G10").Select()
